# USAW Age Groups: Remove spurious weight classes
# - Consolidate championshipName (column B) labels into broader
#   category names: Youth, Juniors, Masters, Seniors
# - Remove the spurious 40kg bodyweight category entry (column K)
#   from the "Open" (row 24) and "Adaptive" (row 25) age groups

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Female section (rows 2-23) ---
$ws.Range("B2").Value  = "Youth"    # U11
$ws.Range("B3").Value  = "Youth"    # U13
$ws.Range("B4").Value  = "Youth"    # U15
$ws.Range("B5").Value  = "Youth"    # U17
$ws.Range("B6").Value  = "Juniors" # U23
$ws.Range("B7").Value  = "Juniors" # U25
$ws.Range("B8").Value  = "Masters" # W35
$ws.Range("B9").Value  = "Masters" # W40
$ws.Range("B10").Value = "Masters" # W45
$ws.Range("B11").Value = "Masters" # W50
$ws.Range("B12").Value = "Masters" # W55
$ws.Range("B13").Value = "Masters" # W60
$ws.Range("B14").Value = "Masters" # W65
$ws.Range("B15").Value = "Masters" # W70
$ws.Range("B16").Value = "Masters" # W75
$ws.Range("B17").Value = "Masters" # W80
$ws.Range("B18").Value = "Masters" # W85
$ws.Range("B19").Value = "Masters" # W90
$ws.Range("B20").Value = "Masters" # W95
$ws.Range("B21").Value = "Masters" # W100
$ws.Range("B22").Value = "Juniors" # JR
$ws.Range("B23").Value = "Seniors" # SR

# Remove spurious 40kg weight class from Open (row 24) and Adaptive (row 25)
$ws.Range("K24").ClearContents()
$ws.Range("K25").ClearContents()

# --- Male section (rows 26-47), same remapping ---
$ws.Range("B26").Value = "Youth"    # U11
$ws.Range("B27").Value = "Youth"    # U13
$ws.Range("B28").Value = "Youth"    # U15
$ws.Range("B29").Value = "Youth"    # U17
$ws.Range("B30").Value = "Juniors" # U23
$ws.Range("B31").Value = "Juniors" # U25
$ws.Range("B32").Value = "Masters" # M35
$ws.Range("B33").Value = "Masters" # M40
$ws.Range("B34").Value = "Masters" # M45
$ws.Range("B35").Value = "Masters" # M50
$ws.Range("B36").Value = "Masters" # M55
$ws.Range("B37").Value = "Masters" # M60
$ws.Range("B38").Value = "Masters" # M65
$ws.Range("B39").Value = "Masters" # M70
$ws.Range("B40").Value = "Masters" # M75
$ws.Range("B41").Value = "Masters" # M80
$ws.Range("B42").Value = "Masters" # M85
$ws.Range("B43").Value = "Masters" # M90
$ws.Range("B44").Value = "Masters" # M95
$ws.Range("B45").Value = "Masters" # M100
$ws.Range("B46").Value = "Juniors" # JR
$ws.Range("B47").Value = "Seniors" # SR
